$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 260.47002222422134
$ws.Range("C2").Value = 290.73128431545058
$ws.Range("D2").Value = 258.88015209052207
$ws.Range("E2").Value = 292.57039203923932

$ws.Range("B3").Value = 261.16854873030132
$ws.Range("C3").Value = 290.36940500889131
$ws.Range("D3").Value = 255.34581346004452
$ws.Range("E3").Value = 299.48533619312377

$ws.Range("B1:E3").Select()
